# Fruta / hortaliza, semanal
# Insert a new weekly record as row 10, pushing the existing rows 10-15 down to 11-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 10 - this shifts rows 10:15 down to 11:16
# (preserving all their existing data/formatting).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44580
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103002
$ws.Range("J10").Value = "Ciruela"
$ws.Range("K10").Value = "Black Amber"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 270
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 1083
$ws.Range("T10").Value = 18
